$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = "'291.83"
$ws.Range('E2').Formula = "'-7.44%"
$ws.Range('G2').Formula = "'6"
$ws.Range('D3').Formula = "'40.39"
$ws.Range('E3').Formula = "'-1.77%"
$ws.Range('G3').Formula = "'6"
$ws.Range('D4').Formula = "'5.032"
$ws.Range('E4').Formula = "'-2.40%"
$ws.Range('G4').Formula = "'6"
$ws.Range('D5').Formula = "'0.07327"
$ws.Range('E5').Formula = "'-3.56%"
$ws.Range('G5').Formula = "'6"
$ws.Range('B6').Value = 'FTXToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D6').Formula = "'1.534"
$ws.Range('E6').Formula = "'-8.28%"
$ws.Range('G6').Formula = "'6"
$ws.Range('B7').Value = 'MXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D7').Formula = "'0.9273"
$ws.Range('E7').Formula = "'-0.23%"
$ws.Range('G7').Formula = "'6"
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').Formula = "'2.348"
$ws.Range('E8').Formula = "'-3.14%"
$ws.Range('G8').Formula = "'6"
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').Formula = "'0.1174"
$ws.Range('E9').Formula = "'-1.94%"
$ws.Range('G9').Formula = "'6"
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Formula = "'0.1746"
$ws.Range('E10').Formula = "'-3.94%"
$ws.Range('G10').Formula = "'6"
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D11').Formula = "'0.04335"
$ws.Range('E11').Formula = "'4.57%"
$ws.Range('G11').Formula = "'6"
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Formula = "'0.08655"
$ws.Range('E12').Formula = "'-4.68%"
$ws.Range('G12').Formula = "'6"
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Formula = "'0.1056"
$ws.Range('E13').Formula = "'0.34%"
$ws.Range('G13').Formula = "'6"
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').Formula = "'0.001268"
$ws.Range('E14').Formula = "'-0.84%"
$ws.Range('G14').Formula = "'6"
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').Formula = "'0.006031"
$ws.Range('E15').Formula = "'4.02%"
$ws.Range('G15').Formula = "'6"
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').Formula = "'3.338"
$ws.Range('E16').Formula = "'0.13%"
$ws.Range('G16').Formula = "'6"
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').Formula = "'4.296"
$ws.Range('E17').Formula = "'-0.67%"
$ws.Range('G17').Formula = "'6"
$ws.Range('E18').Formula = "'-2.01%"
$ws.Range('G18').Formula = "'6"
$ws.Range('D19').Formula = "'7.972"
$ws.Range('E19').Formula = "'4.64%"
$ws.Range('G19').Formula = "'6"
$ws.Range('E20').Formula = "'4.33%"
$ws.Range('G20').Formula = "'6"
$ws.Range('D21').Formula = "'0.2743"
$ws.Range('E21').Formula = "'-7.07%"
$ws.Range('G21').Formula = "'6"
$ws.Range('D22').Formula = "'0.03931"
$ws.Range('E22').Formula = "'-2.30%"
$ws.Range('G22').Formula = "'6"
$ws.Range('E23').Formula = "'-1.10%"
$ws.Range('G23').Formula = "'6"
$ws.Range('D24').Formula = "'0.003779"
$ws.Range('E24').Formula = "'-6.77%"
$ws.Range('G24').Formula = "'6"
$ws.Range('E25').Formula = "'0.81%"
$ws.Range('G25').Formula = "'6"
$ws.Range('D26').Formula = "'0.0003727"
$ws.Range('E26').Formula = "'-95.05%"
$ws.Range('G26').Formula = "'6"
$ws.Range('G27').Formula = "'6"
$ws.Range('G28').Formula = "'6"
$ws.Range('G29').Formula = "'6"
$ws.Range('G30').Formula = "'6"
$ws.Range('G31').Formula = "'6"
$ws.Range('G32').Formula = "'6"
$ws.Range('G33').Formula = "'6"
$ws.Range('G34').Formula = "'6"
$ws.Range('G35').Formula = "'6"
$ws.Range('G36').Formula = "'6"
$ws.Range('G37').Formula = "'6"
$ws.Range('D38').Formula = "'0.02276"
$ws.Range('E38').Formula = "'-5.53%"
$ws.Range('G38').Formula = "'6"
$ws.Range('D39').Formula = "'0.05000"
$ws.Range('E39').Formula = "'-2.96%"
$ws.Range('G39').Formula = "'6"
$ws.Range('D40').Formula = "'0.005911"
$ws.Range('E40').Formula = "'78.97%"
$ws.Range('G40').Formula = "'6"
$ws.Range('D41').Formula = "'0.007686"
$ws.Range('E41').Formula = "'-0.56%"
$ws.Range('G41').Formula = "'6"
$ws.Range('D42').Formula = "'0.1286"
$ws.Range('E42').Formula = "'-0.86%"
$ws.Range('G42').Formula = "'6"
$ws.Range('D43').Formula = "'0.007355"
$ws.Range('E43').Formula = "'-3.13%"
$ws.Range('G43').Formula = "'6"
$ws.Range('D44').Formula = "'0.008302"
$ws.Range('E44').Formula = "'-3.39%"
$ws.Range('G44').Formula = "'6"
$ws.Range('D45').Formula = "'0.2915"
$ws.Range('E45').Formula = "'-14.52%"
$ws.Range('G45').Formula = "'6"
$ws.Range('D46').Formula = "'0.00006282"
$ws.Range('E46').Formula = "'-4.60%"
$ws.Range('G46').Formula = "'6"
$ws.Range('E47').Formula = "'0.02%"
$ws.Range('G47').Formula = "'6"
$ws.Range('D48').Formula = "'0.02837"
$ws.Range('E48').Formula = "'-89.44%"
$ws.Range('G48').Formula = "'6"
$ws.Range('D49').Formula = "'0.00002102"
$ws.Range('E49').Formula = "'0.02%"
$ws.Range('G49').Formula = "'6"
$ws.Range('D50').Formula = "'0.0002002"
$ws.Range('E50').Formula = "'0.02%"
$ws.Range('G50').Formula = "'6"
$ws.Range('G51').Formula = "'6"
